# Auto-generated edit script: apply scheduled-runner price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13 (ALC)
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 28 (ALC)
$ws.Range("H28").Value = 10818.454
$ws.Range("I28").Value = 880
$ws.Range("K28").Value = 880
$ws.Range("M28").Value = -395

# Row 31 (ALC)
$ws.Range("H31").Value = 19.25
$ws.Range("I31").Value = 19.25
$ws.Range("K31").Value = 57.75
$ws.Range("M31").Value = 172.25

# Row 125 (ALC)
$ws.Range("H125").Value = 3436.611
$ws.Range("J125").Value = 3397.8333
$ws.Range("L125").Value = 30580.4997
$ws.Range("N125").Value = -35500.4997

# Row 137 (ALC)
$ws.Range("H137").Value = 2374.739
$ws.Range("I137").Value = 1502
$ws.Range("K137").Value = 4506
$ws.Range("M137").Value = -1956

# Row 138 (ALC)
$ws.Range("H138").Value = 2480.8096
$ws.Range("I138").Value = 907.1667
$ws.Range("K138").Value = 2721.5001
$ws.Range("M138").Value = 2418.4999

$ws = $wb.Worksheets.Item("ARM")
# Row 14 (ARM)
$ws.Range("H14").Value = 6000
$ws.Range("I14").Value = 6000
$ws.Range("K14").Value = 6000
$ws.Range("M14").Value = -5825

# Row 74 (ARM)
$ws.Range("H74").Value = 2991.4167
$ws.Range("I74").Value = 2920.9048
$ws.Range("J74").Value = 3485
$ws.Range("K74").Value = 2920.9048
$ws.Range("L74").Value = 3485
$ws.Range("M74").Value = -2046.9048
$ws.Range("N74").Value = -5233

# Row 77 (ARM)
$ws.Range("H77").Value = 2991.4167
$ws.Range("I77").Value = 2920.9048
$ws.Range("J77").Value = 3485
$ws.Range("K77").Value = 14604.524
$ws.Range("L77").Value = 17425
$ws.Range("M77").Value = -10236.524
$ws.Range("N77").Value = -26161

# Row 132 (ARM)
$ws.Range("H132").Value = 2892.9285
$ws.Range("I132").Value = 2874.5
$ws.Range("K132").Value = 8623.5
$ws.Range("M132").Value = -6093.5

$ws = $wb.Worksheets.Item("BSM")
# Row 14 (BSM)
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

# Row 86 (BSM)
$ws.Range("H86").Value = 2529.9167
$ws.Range("I86").Value = 2529.9167
$ws.Range("K86").Value = 2529.9167
$ws.Range("M86").Value = -1406.9167

# Row 89 (BSM)
$ws.Range("H89").Value = 2529.9167
$ws.Range("I89").Value = 2529.9167
$ws.Range("K89").Value = 12649.5835
$ws.Range("M89").Value = -7033.583500000001

$ws = $wb.Worksheets.Item("CRP")
# Row 2 (CRP)
$ws.Range("H2").Value = 700
$ws.Range("I2").Value = 700
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 700
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -587
$ws.Range("N2").ClearContents()

# Row 31 (CRP)
$ws.Range("H31").Value = 5372.6
$ws.Range("I31").Value = 3957.625
$ws.Range("K31").Value = 3957.625
$ws.Range("M31").Value = -3662.625

# Row 34 (CRP)
$ws.Range("H34").Value = 5372.6
$ws.Range("I34").Value = 3957.625
$ws.Range("K34").Value = 3957.625
$ws.Range("M34").Value = -3755.625

# Row 58 (CRP)
$ws.Range("H58").Value = 4120.75
$ws.Range("I58").Value = 2495.6667
$ws.Range("K58").Value = 2495.6667
$ws.Range("M58").Value = -2292.6667

# Row 99 (CRP)
$ws.Range("H99").Value = 2892
$ws.Range("I99").Value = 2900
$ws.Range("J99").Value = 2882
$ws.Range("K99").Value = 2900
$ws.Range("L99").Value = 2882
$ws.Range("M99").Value = -1402
$ws.Range("N99").Value = -5878

# Row 105 (CRP)
$ws.Range("H105").Value = 1236.6
$ws.Range("I105").Value = 796
$ws.Range("K105").Value = 796
$ws.Range("M105").Value = 951

# Row 126 (CRP)
$ws.Range("H126").Value = 2892
$ws.Range("I126").Value = 2900
$ws.Range("J126").Value = 2882
$ws.Range("K126").Value = 8700
$ws.Range("L126").Value = 8646
$ws.Range("M126").Value = -6230
$ws.Range("N126").Value = -13586

# Row 134 (CRP)
$ws.Range("H134").Value = 2034.591
$ws.Range("I134").Value = 1315.6316
$ws.Range("J134").Value = 6588
$ws.Range("K134").Value = 3946.8948
$ws.Range("L134").Value = 19764
$ws.Range("M134").Value = -1411.8948
$ws.Range("N134").Value = -24834

# Row 136 (CRP)
$ws.Range("H136").Value = 4120.75
$ws.Range("I136").Value = 2495.6667
$ws.Range("K136").Value = 7487.000100000001
$ws.Range("M136").Value = -4937.000100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 6 (CUL)
$ws.Range("H6").Value = 377.35715
$ws.Range("I6").Value = 88.3
$ws.Range("J6").Value = 1100
$ws.Range("K6").Value = 264.9
$ws.Range("L6").Value = 3300
$ws.Range("M6").Value = -151.9
$ws.Range("N6").Value = -3526

# Row 19 (CUL)
$ws.Range("H19").Value = 4999.857
$ws.Range("J19").Value = 4999.857
$ws.Range("L19").Value = 14999.571
$ws.Range("N19").Value = -15347.571

# Row 128 (CUL)
$ws.Range("H128").Value = 484995.16
$ws.Range("I128").Value = 484995.16
$ws.Range("K128").Value = 1454985.48
$ws.Range("M128").Value = -1450005.48

# Row 138 (CUL)
$ws.Range("H138").Value = 4989.857
$ws.Range("J138").Value = 8266.666999999999
$ws.Range("L138").Value = 24800.001
$ws.Range("N138").Value = -35080.001

$ws = $wb.Worksheets.Item("GSM")
# Row 12 (GSM)
$ws.Range("H12").Value = 3486.0667
$ws.Range("I12").Value = 3486.0667
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 3486.0667
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -3346.0667
$ws.Range("N12").ClearContents()

# Row 101 (GSM)
$ws.Range("H101").Value = 14197.6
$ws.Range("J101").Value = 14197.6
$ws.Range("L101").Value = 14197.6
$ws.Range("N101").Value = -20687.6

# Row 104 (GSM)
$ws.Range("H104").Value = 4111.375
$ws.Range("J104").Value = 4111.375
$ws.Range("L104").Value = 4111.375
$ws.Range("N104").Value = -11099.375

# Row 132 (GSM)
$ws.Range("H132").Value = 96672.55
$ws.Range("I132").Value = 130524.75
$ws.Range("K132").Value = 391574.25
$ws.Range("M132").Value = -389044.25

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (LTW)
$ws.Range("H40").Value = 4160.6
$ws.Range("I40").Value = 3034.6667
$ws.Range("K40").Value = 3034.6667
$ws.Range("M40").Value = -2898.6667

# Row 58 (LTW)
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (WVR)
$ws.Range("H122").Value = 2498.84
$ws.Range("I122").Value = 2092.3
$ws.Range("J122").Value = 4125
$ws.Range("K122").Value = 6276.900000000001
$ws.Range("L122").Value = 12375
$ws.Range("M122").Value = -3826.900000000001
$ws.Range("N122").Value = -17275

# Row 132 (WVR)
$ws.Range("H132").Value = 2195.9375
$ws.Range("I132").Value = 2009.6428
$ws.Range("K132").Value = 6028.928400000001
$ws.Range("M132").Value = -3498.928400000001

# Row 136 (WVR)
$ws.Range("H136").Value = 2237.1162
$ws.Range("I136").Value = 1591.0588
$ws.Range("K136").Value = 4773.1764
$ws.Range("M136").Value = -2223.1764
